$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell AB1 with the same formatting as AA1 (border/font/alignment),
# then set its value to the new date label.
$ws.Range("AA1").Copy($ws.Range("AB1"))
$ws.Range("AB1").Value() = "13-10-2020"

# Fill in the new AB column data values (row 2 through row 36).
$ws.Range("AB2").Value() = 3770
$ws.Range("AB3").Value() = 708712
$ws.Range("AB4").Value() = 9403
$ws.Range("AB5").Value() = 166039
$ws.Range("AB6").Value() = 185911
$ws.Range("AB7").Value() = 11898
$ws.Range("AB8").Value() = 116540
$ws.Range("AB9").Value() = 3064
$ws.Range("AB10").Value() = 284844
$ws.Range("AB11").Value() = 33698
$ws.Range("AB12").Value() = 133615
$ws.Range("AB13").Value() = 131228
$ws.Range("AB14").Value() = 14690
$ws.Range("AB15").Value() = 72706
$ws.Range("AB16").Value() = 84461
$ws.Range("AB17").Value() = 592084
$ws.Range("AB18").Value() = 199634
$ws.Range("AB19").Value() = 4126
$ws.Range("AB20").Value() = 130721
$ws.Range("AB21").Value() = 1281896
$ws.Range("AB22").Value() = 10707
$ws.Range("AB23").Value() = 5273
$ws.Range("AB24").Value() = 2046
$ws.Range("AB25").Value() = 5813
$ws.Range("AB26").Value() = 230192
$ws.Range("AB27").Value() = 26555
$ws.Range("AB28").Value() = 112099
$ws.Range("AB29").Value() = 137848
$ws.Range("AB30").Value() = 2925
$ws.Range("AB31").Value() = 607203
$ws.Range("AB32").Value() = 189351
$ws.Range("AB33").Value() = 24623
$ws.Range("AB34").Value() = 47609
$ws.Range("AB35").Value() = 393908
$ws.Range("AB36").Value() = 262103
